$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Target" values (column B) from 60 to 300 for each team row
$ws.Range("B2").Value = 300
$ws.Range("B3").Value = 300
$ws.Range("B4").Value = 300

# Remove the "Shortfall" column (D) - header and values - entirely
$ws.Columns("D").Delete()

# Update selection to match target workbook state
$ws.Range("F6").Select()
